$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.028.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.90%  "
$ws.Range("D3").Value = "'1.827.76"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").Value = "'311.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.98%  "
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").Value = "'0.4645"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.85%  "
$ws.Range("D8").Value = "'0.3750"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.81%  "
$ws.Range("D9").Value = "'0.07283"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.04%  "
$ws.Range("D10").Value = "'0.8656"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.18%  "
$ws.Range("D11").Value = "'19.99"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.31%  "
$ws.Range("E12").Value = "  +7.00%  "
$ws.Range("D13").Value = "'1.848.62"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.78%  "
$ws.Range("D14").Value = "'5.352"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.26%  "
$ws.Range("D15").Value = "'6.550"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.03%  "
$ws.Range("D16").Value = "'91.87"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.02%  "
$ws.Range("D17").Value = "'1.008"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").Value = "'0.000008772"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("D20").Value = "'14.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.37%  "
$ws.Range("D21").Value = "'27.147.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.48%  "
$ws.Range("D22").Value = "'5.165"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.20%  "
$ws.Range("D23").Value = "'10.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.70%  "
$ws.Range("D24").Value = "'2.077.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.36%  "
$ws.Range("D25").Value = "'153.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.13%  "
$ws.Range("D26").Value = "'1.843"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.94%  "
$ws.Range("D27").Value = "'18.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.10%  "
$ws.Range("D28").Value = "'2.087"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.15%  "
$ws.Range("D29").Value = "'5.138"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.57%  "
$ws.Range("D30").Value = "'115.60"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.13%  "
$ws.Range("D31").Value = "'0.08861"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.95%  "
$ws.Range("D32").Value = "'2.965"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.58%  "
$ws.Range("D33").Value = "'0.7308"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.31%  "
$ws.Range("D34").Value = "'4.446"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.01%  "
$ws.Range("E35").Value = "  -3.10%  "
$ws.Range("D36").Value = "'2.483"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.36%  "
$ws.Range("D37").Value = "'1.078"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.46%  "
$ws.Range("D38").Value = "'0.01945"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.37%  "
$ws.Range("D39").Value = "'0.05237"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.72%  "
$ws.Range("D40").Value = "'7.326"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.35%  "
$ws.Range("D41").Value = "'2.922"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.08%  "
$ws.Range("D42").Value = "'0.5169"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.16%  "
$ws.Range("D43").Value = "'0.1630"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.63%  "
$ws.Range("D44").Value = "'0.8578"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -15.05%  "
$ws.Range("D45").Value = "'8.212"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.14%  "
$ws.Range("D46").Value = "'0.4824"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.77%  "
$ws.Range("E47").Value = "  -0.22%  "
$ws.Range("D48").Value = "'10.22"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.99%  "
$ws.Range("D49").Value = "'102.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.98%  "
$ws.Range("D50").Value = "'1.624"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.20%  "
$ws.Range("D51").Value = "'0.06251"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.74%  "
